$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("harp expander v1")

# Update distributor part number for the 100nF capacitor row (new library part)
$ws.Range("H6").Value = "445-6899-2-ND"

# Add new calculation rows below the BOM table
$ws.Range("F34").Formula = "=12/5"
$ws.Range("F35").Formula = "=3/F34"
$ws.Range("F36").Formula = "=1.5/F35"
$ws.Range("F37").Formula = "=1.25/1.5"

# Match style used by the rest of column F (centered, default style)
$ws.Range("F34:F37").HorizontalAlignment = -4108

# Update the sheet selection/active cell to match saved state
$ws.Range("M15").Select() | Out-Null
